$d = $word.ActiveDocument

# --- Add the new closing paragraphs after the table -----------------------
# Paragraph: "AND/OR"
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "AND/OR"

# Blank paragraph
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphAfter()

# Paragraph describing the missing correspondence parts, referencing the
# attachment file instead of the in-document table.
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$r = $p.Range
$r.InsertAfter("Det er ikke likt antall korrespondanseparter som journalposter, da ")

$p = $d.Paragraphs.Last
$r = $p.Range.Duplicate
$r.Collapse(0)
$r.InsertAfter("ANTALLREGISTRERINGERUTENKORRESPONDANSEPART")
$r.Font.Bold = 1
$r.Font.Underline = 1

$p = $d.Paragraphs.Last
$r = $p.Range.Duplicate
$r.Collapse(0)
$r.InsertAfter(" registreringer mangler dette. ")

$p = $d.Paragraphs.Last
$r = $p.Range.Duplicate
$r.Collapse(0)
$r.InsertAfter("Oversikt over j")

$p = $d.Paragraphs.Last
$r = $p.Range.Duplicate
$r.Collapse(0)
$r.InsertAfter("ournalpostene som mangler dette")

$p = $d.Paragraphs.Last
$r = $p.Range.Duplicate
$r.Collapse(0)
$r.InsertAfter(" finnes i vedlegget «3.1.20.txt».")

# --- Tighten the correspondence-parts table columns ------------------------
# (columns 1, 2 and 4 shrink slightly to match the actual cell contents)
$t = $d.Tables.Item(1)
$t.Columns.Item(1).Width = 1642.0 / 20.0
$t.Columns.Item(2).Width = 1166.0 / 20.0
$t.Columns.Item(4).Width = 2343.0 / 20.0
